$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- A1 : "SILAHKAN MASUKKAN`nUSERNAME GURU" (2nd line bold) ---
$a1Text = "SILAHKAN MASUKKAN`nUSERNAME GURU"
$ws.Range("A1").Value = $a1Text
$a1Bold = $ws.Range("A1").Characters(19, 14)
$a1Bold.Font.Bold = $true
$a1Bold.Font.Size = 11
$a1Bold.Font.Name = "Calibri"
$a1Bold.Font.Color = 0

# --- B1 : "SILAHKAN MASUKKAN`nNAMA GURU" (2nd line bold) ---
$b1Text = "SILAHKAN MASUKKAN`nNAMA GURU"
$ws.Range("B1").Value = $b1Text
$b1Bold = $ws.Range("B1").Characters(19, 10)
$b1Bold.Font.Bold = $true
$b1Bold.Font.Size = 11
$b1Bold.Font.Name = "Calibri"
$b1Bold.Font.Color = 0

# --- A1:B1 formatting: fill (accent6 tint), default non-bold font, centered + wrap ---
$hdr = $ws.Range("A1:B1")
$hdr.Interior.ThemeColor = 10
$hdr.Interior.TintAndShade = 0.39997558519241921
$hdr.Font.Bold = $false
$hdr.Font.ThemeColor = 1
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4108
$hdr.WrapText = $true

# --- new row 2 values ---
$ws.Range("A2").Value = "namaguru"
$ws.Range("B2").Value = "Nama Guru"

# --- column widths (closest achievable to 24.77734375 / 23.6640625) ---
$ws.Columns.Item(1).ColumnWidth = 24.0
$ws.Columns.Item(2).ColumnWidth = 22.833333333333332

# --- row 1 custom height ---
$ws.Rows.Item(1).RowHeight = 59.4

# --- selection matching final state ---
$ws.Range("D14").Select()
